$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores prices as plain text (e.g. "48.015.49", "1.00", "0.124")
# using a thousands-dot / decimal-dot mixed display format, never real numbers.
# Assigning some of the new values straight to .Value would make Excel
# auto-detect them as numeric (e.g. "1.00" -> 1, "0.124" -> 0.124) and silently
# drop the exact text formatting (trailing zeros, etc). To avoid that we:
#   1) force the cell to Text format right before the write so Excel keeps
#      the literal string,
#   2) then clear the formatting again so the cell is left with the same
#      "no special style" state it had before the edit (matches the source,
#      where these cells carry no s="" attribute at all).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "48.015.49"
$ws.Range("E2").Value = "  +1.40%  "
Set-TextValue $ws.Range("D3") "2.513.23"
$ws.Range("E3").Value = "  +1.07%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue $ws.Range("D5") "323.01"
$ws.Range("E5").Value = "  +0.09%  "
Set-TextValue $ws.Range("D6") "109.50"
$ws.Range("E6").Value = "  +3.34%  "
$ws.Range("E7").Value = "  +0.03%  "
Set-TextValue $ws.Range("D8") "1.00"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +2.54%  "
Set-TextValue $ws.Range("D10") "40.96"
$ws.Range("E10").Value = "  +7.30%  "
$ws.Range("E11").Value = "  +0.62%  "
Set-TextValue $ws.Range("D12") "0.124"
$ws.Range("E12").Value = "  +0.70%  "
Set-TextValue $ws.Range("D13") "18.77"
$ws.Range("E13").Value = "  +1.94%  "
Set-TextValue $ws.Range("D14") "7.28"
$ws.Range("E14").Value = "  +1.33%  "
Set-TextValue $ws.Range("D15") "2.907.71"
$ws.Range("E15").Value = "  +1.04%  "
Set-TextValue $ws.Range("D16") "2.520.56"
$ws.Range("E16").Value = "  +1.34%  "
Set-TextValue $ws.Range("D17") "0.858"
$ws.Range("E17").Value = "  +1.43%  "
Set-TextValue $ws.Range("D18") "47.903.38"
$ws.Range("E18").Value = "  +1.37%  "
Set-TextValue $ws.Range("D19") "13.32"
$ws.Range("E19").Value = "  +4.29%  "
$ws.Range("E20").Value = "  +1.55%  "
Set-TextValue $ws.Range("D21") "2.81"
$ws.Range("E21").Value = "  +15.85%  "
Set-TextValue $ws.Range("D22") "0.0₃0947"
$ws.Range("E22").Value = "  +0.97%  "
Set-TextValue $ws.Range("D23") "70.97"
$ws.Range("E23").Value = "  +0.49%  "
Set-TextValue $ws.Range("D24") "248.50"
$ws.Range("E24").Value = "  -1.23%  "
Set-TextValue $ws.Range("D25") "2.56"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  -0.47%  "
Set-TextValue $ws.Range("D28") "10.05"
$ws.Range("E28").Value = "  +0.32%  "
Set-TextValue $ws.Range("D29") "2.20"
$ws.Range("E29").Value = "  -0.62%  "
$ws.Range("E30").Value = "  +3.57%  "
Set-TextValue $ws.Range("D31") "35.13"
$ws.Range("E31").Value = "  -0.11%  "
Set-TextValue $ws.Range("D32") "49.76"
$ws.Range("E32").Value = "  +0.57%  "
Set-TextValue $ws.Range("D33") "20.15"
$ws.Range("E33").Value = "  +2.49%  "
Set-TextValue $ws.Range("D34") "5.39"
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("E37").Value = "  +0.46%  "
Set-TextValue $ws.Range("D38") "4.71"
$ws.Range("E38").Value = "  +1.88%  "
Set-TextValue $ws.Range("D39") "2.99"
$ws.Range("E39").Value = "  +0.38%  "
Set-TextValue $ws.Range("D40") "22.64"
$ws.Range("E40").Value = "  +7.39%  "
$ws.Range("E41").Value = "  +0.22%  "
Set-TextValue $ws.Range("D43") "119.58"
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("E44").Value = "  +0.83%  "
Set-TextValue $ws.Range("D45") "2.004.78"
$ws.Range("E45").Value = "  +2.15%  "
Set-TextValue $ws.Range("D46") "3.08"
$ws.Range("E46").Value = "  +3.48%  "
Set-TextValue $ws.Range("D47") "2.03"
$ws.Range("E47").Value = "  -3.27%  "
Set-TextValue $ws.Range("D48") "1.84"
$ws.Range("E48").Value = "  +2.15%  "
Set-TextValue $ws.Range("D49") "9.07"
$ws.Range("E49").Value = "  -0.58%  "
Set-TextValue $ws.Range("D50") "5.23"
$ws.Range("E50").Value = "  -0.60%  "
Set-TextValue $ws.Range("D51") "57.16"
$ws.Range("E51").Value = "  +4.44%  "
